$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws1.Range("B3").Value = -19.364946692964853
$ws1.Range("D3").Value = 26.523555712032437
$ws1.Range("F3").Value = 96.300037659891615
$ws1.Range("G3").Value = 0.017952714771327243
$ws1.Range("H3").Value = 0.84980120462270137
$ws1.Range("I3").Value = 0.54459077402851874

# Row 4 updates
$ws1.Range("B4").Value = -78.496639998254906
$ws1.Range("D4").Value = 24.106221427066124
$ws1.Range("F4").Value = 32.330468954191872
$ws1.Range("G4").Value = 0.06403038201997592
$ws1.Range("H4").Value = 0.97731086153151481
$ws1.Range("I4").Value = 0.60931132471887395
